$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 (pushes "TextBlob" and everything below
# down by one row), matching the workbook's existing sort order
# (Stopwords ISO < SudachiPy < TextBlob alphabetically).
$ws.Rows.Item(33).Insert()

$ws.Range("A33").Value = "SudachiPy"
$ws.Range("B33").Value = "https://github.com/WorksApplications/sudachi.rs"
$ws.Range("C33").Value = "0.6.2"
$ws.Range("D33").Value = "Works Applications Co., Ltd."
$ws.Range("E33").Value = "Apache-2.0"
$ws.Range("F33").Value = "https://github.com/WorksApplications/sudachi.rs/blob/develop/LICENSE"

# Match formatting of the surrounding data rows (A/C/D/E use style index 7,
# B/F use the "hyperlink-look" style index 8).
$ws.Range("A33").Style = $ws.Range("A34").Style
$ws.Range("B33").Style = $ws.Range("B34").Style
$ws.Range("C33").Style = $ws.Range("C34").Style
$ws.Range("D33").Style = $ws.Range("D34").Style
$ws.Range("E33").Style = $ws.Range("E34").Style
$ws.Range("F33").Style = $ws.Range("F34").Style

# Refresh the sort-state range and selection to include the new row.
$ws.Range("C33").Select()
